# chore: update Sheets via scheduled runner
# Applies refreshed market-price values (currentAveragePrice / LevePrice / LeveProfit
# columns H..N) across the ALC, ARM, BSM, CRP, CUL, GSM, LTW and WVR sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# ALC
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ALC")

$ws.Range("H33").Value = 117.89474
$ws.Range("I33").Value = 45
$ws.Range("K33").Value = 45
$ws.Range("M33").Value = 184

$ws.Range("H62").Value = 2902.3076
$ws.Range("I62").Value = 1775.7142
$ws.Range("J62").Value = 4216.6665
$ws.Range("K62").Value = 1775.7142
$ws.Range("L62").Value = 4216.6665
$ws.Range("M62").Value = -1151.7142
$ws.Range("N62").Value = -5464.6665

$ws.Range("H65").Value = 2902.3076
$ws.Range("I65").Value = 1775.7142
$ws.Range("J65").Value = 4216.6665
$ws.Range("K65").Value = 8878.571
$ws.Range("L65").Value = 21083.3325
$ws.Range("M65").Value = -5758.571
$ws.Range("N65").Value = -27323.3325

$ws.Range("H125").Value = 2537.353
$ws.Range("I125").Value = 1642.7142
$ws.Range("J125").Value = 3163.6
$ws.Range("K125").Value = 14784.4278
$ws.Range("L125").Value = 28472.4
$ws.Range("M125").Value = -12324.4278
$ws.Range("N125").Value = -33392.39999999999

$ws.Range("H137").Value = 3201
$ws.Range("I137").Value = 2011.7693
$ws.Range("K137").Value = 6035.3079
$ws.Range("M137").Value = -3485.3079

$ws.Range("H138").Value = 2448.58
$ws.Range("I138").Value = 786.7742
$ws.Range("K138").Value = 2360.3226
$ws.Range("M138").Value = 2779.6774

$ws.Range("H141").Value = 2861.4395
$ws.Range("I141").Value = 2653.125
$ws.Range("J141").Value = 4028
$ws.Range("K141").Value = 7959.375
$ws.Range("L141").Value = 12084
$ws.Range("M141").Value = -2779.375
$ws.Range("N141").Value = -22444

# ---------------------------------------------------------------------------
# ARM
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ARM")

$ws.Range("H32").Value = 2997.602
$ws.Range("I32").Value = 2744.5784
$ws.Range("J32").Value = 5097.7
$ws.Range("K32").Value = 2744.5784
$ws.Range("L32").Value = 5097.7
$ws.Range("M32").Value = -2457.5784
$ws.Range("N32").Value = -5671.7

$ws.Range("H61").Value = 1106.52
$ws.Range("I61").Value = 811.6111
$ws.Range("J61").Value = 1864.8572
$ws.Range("K61").Value = 811.6111
$ws.Range("L61").Value = 1864.8572
$ws.Range("M61").Value = -599.6111
$ws.Range("N61").Value = -2288.8572

$ws.Range("H74").Value = 2242.6128
$ws.Range("I74").Value = 2050.6072
$ws.Range("J74").Value = 4034.6667
$ws.Range("K74").Value = 2050.6072
$ws.Range("L74").Value = 4034.6667
$ws.Range("M74").Value = -1176.6072
$ws.Range("N74").Value = -5782.6667

$ws.Range("H77").Value = 2242.6128
$ws.Range("I77").Value = 2050.6072
$ws.Range("J77").Value = 4034.6667
$ws.Range("K77").Value = 10253.036
$ws.Range("L77").Value = 20173.3335
$ws.Range("M77").Value = -5885.036
$ws.Range("N77").Value = -28909.3335

$ws.Range("H132").Value = 2631.475
$ws.Range("I132").Value = 1781.9615
$ws.Range("J132").Value = 4209.143
$ws.Range("K132").Value = 5345.8845
$ws.Range("L132").Value = 12627.429
$ws.Range("M132").Value = -2815.8845
$ws.Range("N132").Value = -17687.429

$ws.Range("H136").Value = 1106.52
$ws.Range("I136").Value = 811.6111
$ws.Range("J136").Value = 1864.8572
$ws.Range("K136").Value = 2434.8333
$ws.Range("L136").Value = 5594.571599999999
$ws.Range("M136").Value = 115.1667000000002
$ws.Range("N136").Value = -10694.5716

# ---------------------------------------------------------------------------
# BSM
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("BSM")

$ws.Range("H94").Value = 453.625
$ws.Range("I94").Value = 557.2273
$ws.Range("J94").Value = 225.7
$ws.Range("K94").Value = 557.2273
$ws.Range("L94").Value = 225.7
$ws.Range("M94").Value = -106.2273
$ws.Range("N94").Value = -1127.7

$ws.Range("H134").Value = 1892.0476
$ws.Range("I134").Value = 1050.2157
$ws.Range("J134").Value = 5469.8335
$ws.Range("K134").Value = 3150.6471
$ws.Range("L134").Value = 16409.5005
$ws.Range("M134").Value = -615.6471000000001
$ws.Range("N134").Value = -21479.5005

# ---------------------------------------------------------------------------
# CRP
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("CRP")

$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("M2").ClearContents()

$ws.Range("H31").Value = 8477177
$ws.Range("I31").Value = 1364.8064
$ws.Range("J31").Value = 17861112
$ws.Range("K31").Value = 1364.8064
$ws.Range("L31").Value = 17861112
$ws.Range("M31").Value = -1069.8064
$ws.Range("N31").Value = -17861702

$ws.Range("H34").Value = 8477177
$ws.Range("I34").Value = 1364.8064
$ws.Range("J34").Value = 17861112
$ws.Range("K34").Value = 1364.8064
$ws.Range("L34").Value = 17861112
$ws.Range("M34").Value = -1162.8064
$ws.Range("N34").Value = -17861516

$ws.Range("H58").Value = 1229.06
$ws.Range("I58").Value = 1355.117
$ws.Range("J58").Value = 807.04346
$ws.Range("K58").Value = 1355.117
$ws.Range("L58").Value = 807.04346
$ws.Range("M58").Value = -1152.117
$ws.Range("N58").Value = -1213.04346

$ws.Range("H99").Value = 9096147
$ws.Range("I99").Value = 14289380
$ws.Range("J99").Value = 7989.25
$ws.Range("K99").Value = 14289380
$ws.Range("L99").Value = 7989.25
$ws.Range("M99").Value = -14287882
$ws.Range("N99").Value = -10985.25

$ws.Range("H126").Value = 9096147
$ws.Range("I126").Value = 14289380
$ws.Range("J126").Value = 7989.25
$ws.Range("K126").Value = 42868140
$ws.Range("L126").Value = 23967.75
$ws.Range("M126").Value = -42865670
$ws.Range("N126").Value = -28907.75

$ws.Range("H132").Value = 5694.154
$ws.Range("I132").Value = 4929.4546
$ws.Range("K132").Value = 14788.3638
$ws.Range("M132").Value = -12258.3638

$ws.Range("H136").Value = 1229.06
$ws.Range("I136").Value = 1355.117
$ws.Range("J136").Value = 807.04346
$ws.Range("K136").Value = 4065.351
$ws.Range("L136").Value = 2421.13038
$ws.Range("M136").Value = -1515.351
$ws.Range("N136").Value = -7521.130380000001

# ---------------------------------------------------------------------------
# CUL
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("CUL")

$ws.Range("H5").Value = 1328.238
$ws.Range("I5").Value = 325.21054
$ws.Range("K5").Value = 975.6316199999999
$ws.Range("M5").Value = -863.6316199999999

$ws.Range("H8").Value = 93.25
$ws.Range("I8").Value = 93.25
$ws.Range("K8").Value = 279.75
$ws.Range("M8").Value = -140.75

$ws.Range("H12").Value = 125
$ws.Range("I12").Value = 10.5
$ws.Range("J12").Value = 150.44444
$ws.Range("K12").Value = 31.5
$ws.Range("L12").Value = 451.33332
$ws.Range("M12").Value = 141.5
$ws.Range("N12").Value = -797.33332

$ws.Range("H87").Value = 3014
$ws.Range("I87").Value = 3014
$ws.Range("J87").Value = 0
$ws.Range("K87").Value = 9042
$ws.Range("L87").Value = 0
$ws.Range("M87").Value = -7794
$ws.Range("N87").ClearContents()

$ws.Range("H90").Value = 3014
$ws.Range("I90").Value = 3014
$ws.Range("J90").Value = 0
$ws.Range("K90").Value = 27126
$ws.Range("L90").Value = 0
$ws.Range("M90").Value = -20886
$ws.Range("N90").ClearContents()

$ws.Range("H113").Value = 532.06665
$ws.Range("I113").Value = 487.51428
$ws.Range("K113").Value = 1462.54284
$ws.Range("M113").Value = 707.4571599999999

$ws.Range("H135").Value = 1328.238
$ws.Range("I135").Value = 325.21054
$ws.Range("K135").Value = 2926.89486
$ws.Range("M135").Value = -391.8948599999999

$ws.Range("H137").Value = 2579.7646
$ws.Range("I137").Value = 877.8
$ws.Range("J137").Value = 5011.143
$ws.Range("K137").Value = 2633.4
$ws.Range("L137").Value = 15033.429
$ws.Range("M137").Value = 2466.6
$ws.Range("N137").Value = -25233.429

$ws.Range("H138").Value = 2953
$ws.Range("I138").Value = 2604
$ws.Range("K138").Value = 7812
$ws.Range("M138").Value = -2672

# ---------------------------------------------------------------------------
# GSM
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("GSM")

$ws.Range("H70").Value = 6496.1143
$ws.Range("I70").Value = 5865.609
$ws.Range("K70").Value = 5865.609
$ws.Range("M70").Value = -5595.609

$ws.Range("H73").Value = 6496.1143
$ws.Range("I73").Value = 5865.609
$ws.Range("K73").Value = 5865.609
$ws.Range("M73").Value = -4929.609

$ws.Range("H132").Value = 3013.1785
$ws.Range("I132").Value = 1983.4706
$ws.Range("J132").Value = 4604.5454
$ws.Range("K132").Value = 5950.4118
$ws.Range("L132").Value = 13813.6362
$ws.Range("M132").Value = -3420.4118
$ws.Range("N132").Value = -18873.6362

# ---------------------------------------------------------------------------
# LTW
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("LTW")

$ws.Range("H40").Value = 40006856
$ws.Range("I40").Value = 83337120
$ws.Range("K40").Value = 83337120
$ws.Range("M40").Value = -83336984

$ws.Range("H61").Value = 1330.7
$ws.Range("I61").Value = 1122.4286
$ws.Range("J61").Value = 1816.6666
$ws.Range("K61").Value = 1122.4286
$ws.Range("L61").Value = 1816.6666
$ws.Range("M61").Value = -920.4286
$ws.Range("N61").Value = -2220.6666

$ws.Range("H68").Value = 747.4949
$ws.Range("I68").Value = 641.56525
$ws.Range("K68").Value = 641.56525
$ws.Range("M68").Value = 107.43475

$ws.Range("H71").Value = 747.4949
$ws.Range("I71").Value = 641.56525
$ws.Range("K71").Value = 3207.82625
$ws.Range("M71").Value = 536.1737499999999

$ws.Range("H113").Value = 1330.7
$ws.Range("I113").Value = 1122.4286
$ws.Range("J113").Value = 1816.6666
$ws.Range("K113").Value = 1122.4286
$ws.Range("L113").Value = 1816.6666
$ws.Range("M113").Value = 1047.5714
$ws.Range("N113").Value = -6156.6666

$ws.Range("H122").Value = 7023.25
$ws.Range("I122").Value = 3318.5
$ws.Range("K122").Value = 9955.5
$ws.Range("M122").Value = -7505.5

$ws.Range("H132").Value = 17468.727
$ws.Range("I132").Value = 19900.857
$ws.Range("J132").Value = 13212.5
$ws.Range("K132").Value = 59702.571
$ws.Range("L132").Value = 39637.5
$ws.Range("M132").Value = -57172.571
$ws.Range("N132").Value = -44697.5

$ws.Range("H136").Value = 2982.0557
$ws.Range("I136").Value = 1759.7693
$ws.Range("K136").Value = 5279.3079
$ws.Range("M136").Value = -2729.3079

# ---------------------------------------------------------------------------
# WVR
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("WVR")

$ws.Range("H122").Value = 3267.7307
$ws.Range("I122").Value = 2165.3333
$ws.Range("J122").Value = 5748.125
$ws.Range("K122").Value = 6495.999899999999
$ws.Range("L122").Value = 17244.375
$ws.Range("M122").Value = -4045.999899999999
$ws.Range("N122").Value = -22144.375

$ws.Range("H132").Value = 18520862
$ws.Range("I132").Value = 989.2727
$ws.Range("K132").Value = 2967.8181
$ws.Range("M132").Value = -437.8181

$ws.Range("H136").Value = 1514.4133
$ws.Range("I136").Value = 415.57626
$ws.Range("J136").Value = 5566.375
$ws.Range("K136").Value = 1246.72878
$ws.Range("L136").Value = 16699.125
$ws.Range("M136").Value = 1303.27122
$ws.Range("N136").Value = -21799.125
